# ADD results from server
# Updates computed result values (row 2) on the "2025", "2030" and "2035"
# sheets of the inv_capacity results workbook with refreshed numbers
# produced by the server-side run.

$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 0.002269897435823648
$ws.Range("E2").Value = 0.3191563207764054
$ws.Range("I2").Value = 0.3412040122747214
$ws.Range("L2").Value = 0.5256036900000001
$ws.Range("M2").Value = 0.07430661880348029
$ws.Range("N2").Value = 11.5903253386221
$ws.Range("O2").Value = 3.093707044758613

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0.006578532693934032
$ws.Range("B2").Value = 0.0581238
$ws.Range("E2").Value = 0.3304135492235945
$ws.Range("I2").Value = 0.5856753215803109
$ws.Range("L2").Value = 0.1881759299999999
$ws.Range("M2").Value = 0.07332068119651973
$ws.Range("N2").Value = 8.092351482587798
$ws.Range("O2").Value = 3.347862294350149

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 0.08815316193520684
$ws.Range("B2").Value = 0.02186100000000002
$ws.Range("E2").Value = 0.1406309873331282
$ws.Range("I2").Value = 0.395620046908193
$ws.Range("M2").Value = 0.04415174999999996
$ws.Range("N2").Value = 7.773198087410442
$ws.Range("O2").Value = 4.898127034979207
